$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H1: date value (2015-01-01, serial 42005) formatted as a short date (built-in numFmtId 14)
$ws.Range("H1").Value = 42005
$ws.Range("H1").NumberFormat = "mm-dd-yy"

# H2:H7: a single space character, same shared string as column C uses
$ws.Range("H2").Value = " "
$ws.Range("H3").Value = " "
$ws.Range("H4").Value = " "
$ws.Range("H5").Value = " "
$ws.Range("H6").Value = " "
$ws.Range("H7").Value = " "

# Selection moves to H7
$ws.Range("H7").Select() | Out-Null
